$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

$ws.Range("A$row").Value = 5
$ws.Range("B$row").Value = "Macroferia Regional de Talca"
$ws.Range("C$row").Value = "Maule"
$ws.Range("D$row").Value = 44516
$ws.Range("D$row").NumberFormat = $ws.Range("D24").NumberFormat
$ws.Range("E$row").Value = 7
$ws.Range("F$row").Value = "Fruta"
$ws.Range("G$row").Value = 100101
$ws.Range("H$row").Value = "Berries"
$ws.Range("I$row").Value = 100101001
$ws.Range("J$row").Value = "Arándano (blue)"
$ws.Range("K$row").Value = "Sin especificar"
$ws.Range("L$row").Value = "Primera"
$ws.Range("M$row").Value = 30
$ws.Range("N$row").Value = 5000
$ws.Range("O$row").Value = 5000
$ws.Range("P$row").Value = 5000
$ws.Range("Q$row").Value = "`$/bandeja 2 kilos"
$ws.Range("R$row").Value = "Provincia de Linares"
$ws.Range("S$row").Value = 2500
$ws.Range("T$row").Value = 2
